$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.651.71"
$ws.Range("E2").Value = "  +3.08%  "
$ws.Range("D3").Value = "3.151.60"
$ws.Range("E3").Value = "  +2.60%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'578.61"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").Value = "'180.13"
$ws.Range("E6").Value = "  +6.53%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.152.02"
$ws.Range("E8").Value = "  +2.71%  "
$ws.Range("D9").Value = "'0.524"
$ws.Range("E9").Value = "  +2.72%  "
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("D11").Value = "'0.153"
$ws.Range("E11").Value = "  +2.39%  "
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("E13").Value = "  +2.12%  "
$ws.Range("D14").Value = "'37.10"
$ws.Range("E14").Value = "  +4.38%  "
$ws.Range("D15").Value = "68.644.16"
$ws.Range("E15").Value = "  +3.16%  "
$ws.Range("D16").Value = "3.674.52"
$ws.Range("E16").Value = "  +2.59%  "
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").Value = "'7.17"
$ws.Range("E18").Value = "  +3.38%  "
$ws.Range("D19").Value = "3.148.38"
$ws.Range("E19").Value = "  +2.17%  "
$ws.Range("E20").Value = "  -2.36%  "
$ws.Range("D21").Value = "'489.92"
$ws.Range("E21").Value = "  +0.75%  "
$ws.Range("D22").Value = "'0.701"
$ws.Range("E22").Value = "  +2.30%  "
$ws.Range("E23").Value = "  +1.70%  "
$ws.Range("D24").Value = "'84.10"
$ws.Range("E25").Value = "  +6.99%  "
$ws.Range("E26").Value = "  +3.29%  "
$ws.Range("D27").Value = "'10.63"
$ws.Range("E27").Value = "  +4.94%  "
$ws.Range("D29").Value = "'8.14"
$ws.Range("E29").Value = "  +4.69%  "
$ws.Range("D30").Value = "'2.37"
$ws.Range("E30").Value = "  +4.90%  "
$ws.Range("E31").Value = "  +2.07%  "
$ws.Range("D32").Value = "'28.31"
$ws.Range("E32").Value = "  +3.03%  "
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("D34").Value = "0.0₃0950"
$ws.Range("E34").Value = "  +4.83%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  +3.61%  "
$ws.Range("D37").Value = "'48.28"
$ws.Range("E37").Value = "  +2.58%  "
$ws.Range("D38").Value = "'0.964"
$ws.Range("E38").Value = "  +1.72%  "
$ws.Range("E39").Value = "  +8.72%  "
$ws.Range("D40").Value = "'2.05"
$ws.Range("E40").Value = "  +4.50%  "
$ws.Range("E41").Value = "  +3.11%  "
$ws.Range("D42").Value = "'49.19"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").Value = "'8.43"
$ws.Range("E43").Value = "  +1.92%  "
$ws.Range("E44").Value = "  +9.19%  "
$ws.Range("D45").Value = "'403.65"
$ws.Range("E45").Value = "  +10.44%  "
$ws.Range("D46").Value = "'28.05"
$ws.Range("E46").Value = "  +15.27%  "
$ws.Range("D47").Value = "2.812.62"
$ws.Range("E47").Value = "  +1.72%  "
$ws.Range("E48").Value = "  +1.74%  "
$ws.Range("D49").Value = "'135.30"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("D51").Value = "'2.37"
$ws.Range("E51").Value = "  +10.50%  "
